$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.840.91'
$ws.Range('E2').Value = '  +3.47%  '

# Row 3
$ws.Range('D3').Value = '3.796.02'
$ws.Range('E3').Value = '  +3.95%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.76%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '421.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.06%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.75%  '

# Row 7
$ws.Range('D7').Value = '3.794.40'
$ws.Range('E7').Value = '  +3.98%  '

# Row 8
$ws.Range('E8').Value = '  -2.56%  '

# Row 9
$ws.Range('E9').Value = '  -0.13%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.717'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.36%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.160'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.06%  '

# Row 12
$ws.Range('E12').Value = '  +15.90%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.38'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.29%  '

# Row 14
$ws.Range('D14').Value = '4.409.57'
$ws.Range('E14').Value = '  +4.26%  '

# Row 15
$ws.Range('E15').Value = '  +2.74%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +15.68%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.832.79'
$ws.Range('E17').Value = '  +5.25%  '

# Row 18
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.137'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.57%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.07%  '

# Row 20
$ws.Range('D20').Value = '66.930.20'
$ws.Range('E20').Value = '  +3.18%  '

# Row 21
$ws.Range('E21').Value = '  +0.77%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '403.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.32%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.45'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.67%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.19%  '

# Row 25
$ws.Range('E25').Value = '  +1.68%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.84%  '

# Row 27
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.94%  '

# Row 28
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.66%  '

# Row 29
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.11%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +29.48%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '725.92'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.95%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.50'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.28%  '

# Row 33
$ws.Range('E33').Value = '  +2.50%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.120'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.56%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.15%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.153'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.34%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.74%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.07'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.25%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +24.22%  '

# Row 40
$ws.Range('D40').Value = '0.0₃0751'
$ws.Range('E40').Value = '  +22.29%  '

# Row 41
$ws.Range('E41').Value = '  -2.15%  '

# Row 42
$ws.Range('E42').Value = '  +0.69%  '

# Row 43
$ws.Range('E43').Value = '  +0.45%  '

# Row 44
$ws.Range('E44').Value = '  -3.99%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.09%  '

# Row 46
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '143.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.85%  '

# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.88%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.307'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.75%  '

# Row 49
$ws.Range('E49').Value = '  -1.07%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.23%  '

# Row 51
$ws.Range('E51').Value = '  +1.93%  '
